# Adds four new leetcode entries (with hyperlinks) to the "算法题记录" sheet:
#   F1 -> "环形链表"        (Linked List Cycle)
#   D4 -> "反转链表"        (Reverse Linked List)
#   D3 -> "反转链表"        (Reverse Linked List)
#   E4 -> "K个一组翻转链表"  (Reverse Nodes in k-Group)
# matching the same "text + hyperlink + hyperlink-style" pattern already
# used by the other cells on the sheet (C1:C4, D1, E1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F1: 环形链表 ---------------------------------------------------------
$ws.Range("F1").Value = "环形链表"
$ws.Hyperlinks.Add($ws.Range("F1"), "https://leetcode-cn.com/problems/linked-list-cycle/")
$ws.Range("F1").Style = "超链接"

# --- D4: 反转链表 ---------------------------------------------------------
$ws.Range("D4").Value = "反转链表"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://leetcode-cn.com/problems/reverse-linked-list/")
$ws.Range("D4").Style = "超链接"

# --- D3: 反转链表 ---------------------------------------------------------
$ws.Range("D3").Value = "反转链表"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://leetcode-cn.com/problems/reverse-linked-list/")
$ws.Range("D3").Style = "超链接"

# --- E4: K个一组翻转链表 ---------------------------------------------------
$ws.Range("E4").Value = "K个一组翻转链表"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://leetcode-cn.com/problems/reverse-nodes-in-k-group/")
$ws.Range("E4").Style = "超链接"

# Column E needs to widen to fit the new, longer "K个一组翻转链表" text.
$ws.Columns("E").ColumnWidth = 15.5

# Final selection lands on the last-edited cell, E4.
$ws.Range("E4").Select()
